$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per refreshed crypto feed.
# Cells whose new price text would parse as a plain number (e.g. "211.44")
# are force-formatted as Text first so they round-trip as strings, just like
# the existing sheet data (matches values such as "27.885.07" that keep their
# text type naturally because of the embedded thousand separators).
$ws.Range("D2").Value = "27.885.07"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.630.55"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.44"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.45"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0883"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.863.67"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "1.634.49"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.39"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "27.896.19"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.58"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.05"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.63"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.53"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").Value = "1.393.36"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  +9.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.557"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.848"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "1.773.80"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.14"
$ws.Range("E47").Value = "  -3.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.56"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  +0.84%  "
